# Update "Hoja1" risk-data table:
#  - revise the Dec-2024 (row 13) and Jan-2025 (row 14) figures
#  - add a new Feb-2025 (row 15) observation
#  - push the old blank spacer row down to row 16
#  - move the active selection to G11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- keep formatting consistent while the table grows ---------------------
# The blank spacer that used to live in row 15 (style-only, no data) now
# belongs one row lower, at row 16 - copy its look down first ...
$ws.Range("B15:D15").Copy()
$ws.Range("B16:D16").PasteSpecial(-4122)

# ... then give the brand-new data row 15 the same look as row 14 (date +
# percentage formatting) before filling in its numbers.
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)

# --- revised Dec-2024 figures (row 13) -------------------------------------
$ws.Range("B13").Value = 0.0747308210715293
$ws.Range("C13").Value = 0.046955098168393286
$ws.Range("D13").Value = 0.02309145656447239

# --- revised Jan-2025 figures (row 14) -------------------------------------
$ws.Range("B14").Value = 0.06248937940529241
$ws.Range("C14").Value = 0.06029191320942946
$ws.Range("D14").Value = 0.027301263152702267
$ws.Rows.Item(14).RowHeight = 15

# --- new Feb-2025 observation (row 15) --------------------------------------
$ws.Range("A15").Value = 45701
$ws.Range("B15").Value = 0.04858630462881415
$ws.Range("C15").Value = 0.0641946244709794
$ws.Range("D15").Value = 0.015885432933623746

# --- move the active cell selection -----------------------------------------
$ws.Range("G11").Select()

$excel.CutCopyMode = $false
